$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Innsbruck is a city in wich european country?"
$ws.Range("B2").Value = "Austria"
$ws.Range("C2").Value = "Location"

$ws.Range("A3").Value = "When did Miachel Schumacher win his first F1 World Drivers Title?"
$ws.Range("B3").Value = 1994
$ws.Range("C3").Value = "Year"

$ws.Range("A4").Value = "Who was the F1 World Champion in 2022?"
$ws.Range("B4").Value = "Max Verstappen"
$ws.Range("C4").Value = "Person"

[void]$ws.Range("D10").Select()
